# Generate Report for Handback
# The b889dd7b-7298-4181-abc3-993d682d418a.md file has now been handed back
# (both locales are in sync with en-US), so flip its status from
# "Ready for handoff" to "Handed back: in sync with en-US" and stamp the
# new handback datetimes.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the b889dd7b...md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# --- zh-cn sheet: row 3 is the b889dd7b...md file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusText
$wsZhCn.Range("H3").Value = "2016-03-17 22:38:39"

# --- de-de sheet: row 3 is the b889dd7b...md file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusText
$wsDeDe.Range("H3").Value = "2016-03-17 22:38:45"
